# "Generate Report for Archive"
# The localization-status report is regenerated: the entry for
# cbae1ca0-eb85-4205-a5eb-958f2e3718c3.md moves from "Ready for handoff"
# to "In Translation" on every sheet that tracks its status.

$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn (col E) and de-de (col F) status columns for row 5
# (cbae1ca0-eb85-4205-a5eb-958f2e3718c3.md)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"

# zh-cn sheet: Status column (C) for row 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"

# de-de sheet: Status column (C) for row 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
